# Apply updated values to Sheet2 and Sheet3 after clipping binary masks to [0, 1]
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Sheet2")
$ws.Range("B2").Value = 0.00124669307300026
$ws.Range("C2").Value = 0.02086529610309911
$ws.Range("D2").Value = 0.1611430503889077
$ws.Range("E2").Value = 0
$ws.Range("F2").Value = 0.2860575767770462
$ws.Range("G2").Value = 0.00002933411557641537
$ws.Range("I2").Value = 0.03224457279426409
$ws.Range("B3").Value = 0.03330472564329513
$ws.Range("C3").Value = 0.07755434396501001
$ws.Range("D3").Value = 0.1812411555871863
$ws.Range("E3").Value = 0.01104776235585715
$ws.Range("F3").Value = 0.3355435166647038
$ws.Range("G3").Value = 0.04398504287797503
$ws.Range("H3").Value = 0.1955116511847632
$ws.Range("I3").Value = 0.1838140241453101
$ws.Range("B4").Value = 0.2157046534868357
$ws.Range("C4").Value = 0.08896976409543406
$ws.Range("D4").Value = 0.1775775833621003
$ws.Range("E4").Value = 0.0477085264313921
$ws.Range("F4").Value = 0.647393858279426
$ws.Range("G4").Value = 0.045964898724823
$ws.Range("H4").Value = 0.1284133973040299
$ws.Range("I4").Value = 0.3050628365290755

$ws = $wb.Worksheets.Item("Sheet3")
$ws.Range("B2").Value = 0
$ws.Range("C2").Value = 0
$ws.Range("D2").Value = 0
$ws.Range("E2").Value = 0
$ws.Range("F2").Value = 0
$ws.Range("H2").Value = 1
$ws.Range("I2").Value = 0
$ws.Range("B3").Value = 0.004827801787229398
$ws.Range("C3").Value = 0.01212158443433304
$ws.Range("D3").Value = 0.02105518643245047
$ws.Range("E3").Value = 0.00008471325169396249
$ws.Range("F3").Value = 0
$ws.Range("G3").Value = 0.000560878709978967
$ws.Range("H3").Value = 0.00001379816594778222
$ws.Range("I3").Value = 0.03588458870079315
$ws.Range("B4").Value = 0.1567641988851329
$ws.Range("C4").Value = 0.04361828263447123
$ws.Range("D4").Value = 0.01883848451461288
$ws.Range("E4").Value = 0.02473802668608541
$ws.Range("F4").Value = 0.02280500569374648
$ws.Range("G4").Value = 0.0004942128287888125
$ws.Range("H4").Value = 0.008357239595376457
$ws.Range("I4").Value = 0.07988230396342109
